# [Word] (tracked changes) Clean up and mappings (#891)
#
# Applies the data changes described by the commit to the "Snippets" sheet:
#  - Clears the (already-empty) C14/C15 helper cells (AnnotationPopupActionEventArgs
#    rows) that no longer carry an explicit format.
#  - Re-derives the TrackedChange / TrackedChangeCollection snippet rows
#    (310-317): the "class" header row for TrackedChange now documents
#    getNextTrackedChange instead of acceptFirstTrackedChange, and every
#    member row shifts up to make room for TrackedChangeCollection.rejectAll.
#  - Adds two new snippet rows: the TrackedChangeType enum header row, and
#    the (re-numbered) VerticalAlignment enum header row.
#  - Grows the "Snippets" table/autofilter + sheet dimension to A1:F319.
#  - Restores the sheet view to the top of the table (was scrolled/selected
#    around row 144-160).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Annotation popup rows: drop the stray empty C14 / C15 cells.
# ---------------------------------------------------------------------
$ws.Range("C14").Clear()
$ws.Range("C15").Clear()

# ---------------------------------------------------------------------
# 2. Grow the "Snippets" table to fit the two extra rows (319 total).
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F319"))

# ---------------------------------------------------------------------
# 3. Rewrite the TrackedChange / TrackedChangeCollection block (310-317)
#    plus the two brand-new rows (318-319).
# ---------------------------------------------------------------------

# Row 310: TrackedChange "class" header row now points at getNextTrackedChange.
$ws.Range("A310").Value = "Word"
$ws.Range("B310").Value = "TrackedChange"
$ws.Range("C310").Clear()
$ws.Range("D310").Value = "class"
$ws.Range("E310").Value = "word-document-manage-tracked-changes"
$ws.Range("F310").Value = "getNextTrackedChange"

# Row 311: TrackedChange.accept
$ws.Range("A311").Value = "Word"
$ws.Range("B311").Value = "TrackedChange"
$ws.Range("C311").Value = "accept"
$ws.Range("D311").Value = 1
$ws.Range("E311").Value = "word-document-manage-tracked-changes"
$ws.Range("F311").Value = "acceptFirstTrackedChange"

# Row 312: TrackedChange.getNext
$ws.Range("A312").Value = "Word"
$ws.Range("B312").Value = "TrackedChange"
$ws.Range("C312").Value = "getNext"
$ws.Range("D312").Value = 1
$ws.Range("E312").Value = "word-document-manage-tracked-changes"
$ws.Range("F312").Value = "getNextTrackedChange"

# Row 313: TrackedChange.getRange
$ws.Range("A313").Value = "Word"
$ws.Range("B313").Value = "TrackedChange"
$ws.Range("C313").Value = "getRange"
$ws.Range("D313").Value = 1
$ws.Range("E313").Value = "word-document-manage-tracked-changes"
$ws.Range("F313").Value = "getFirstTrackedChangeRange"

# Row 314: TrackedChange.reject
$ws.Range("A314").Value = "Word"
$ws.Range("B314").Value = "TrackedChange"
$ws.Range("C314").Value = "reject"
$ws.Range("D314").Value = 1
$ws.Range("E314").Value = "word-document-manage-tracked-changes"
$ws.Range("F314").Value = "rejectFirstTrackedChange"

# Row 315: TrackedChangeCollection.acceptAll
$ws.Range("A315").Value = "Word"
$ws.Range("B315").Value = "TrackedChangeCollection"
$ws.Range("C315").Value = "acceptAll"
$ws.Range("D315").Value = 1
$ws.Range("E315").Value = "word-document-manage-tracked-changes"
$ws.Range("F315").Value = "acceptAllTrackedChanges"

# Row 316: TrackedChangeCollection.getFirst
$ws.Range("A316").Value = "Word"
$ws.Range("B316").Value = "TrackedChangeCollection"
$ws.Range("C316").Value = "getFirst"
$ws.Range("D316").Value = 1
$ws.Range("E316").Value = "word-document-manage-tracked-changes"
$ws.Range("F316").Value = "getFirstTrackedChangeRange"

# Row 317 (new): TrackedChangeCollection.rejectAll
$ws.Range("A317").Value = "Word"
$ws.Range("B317").Value = "TrackedChangeCollection"
$ws.Range("C317").Value = "rejectAll"
$ws.Range("D317").Value = 1
$ws.Range("E317").Value = "word-document-manage-tracked-changes"
$ws.Range("F317").Value = "rejectAllTrackedChanges"

# Row 318 (new): TrackedChangeType enum header row
$ws.Range("A318").Value = "Word"
$ws.Range("B318").Value = "TrackedChangeType"
$ws.Range("C318").Clear()
$ws.Range("D318").Value = "enum"
$ws.Range("E318").Value = "word-document-manage-tracked-changes"
$ws.Range("F318").Value = "getNextTrackedChange"

# Row 319 (renumbered from old 317): VerticalAlignment enum header row
$ws.Range("A319").Value = "Word"
$ws.Range("B319").Value = "VerticalAlignment"
$ws.Range("C319").Clear()
$ws.Range("D319").Value = "enum"
$ws.Range("E319").Value = "word-tables-manage-formatting"
$ws.Range("F319").Value = "getTableRowAlignment"

# ---------------------------------------------------------------------
# 4. Reset the view: scroll back to the top of the frozen header and
#    select A2 (was left scrolled to around row 144, cell F160 selected).
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null

Write-Output "Snippets sheet updated through row 319"
